$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.12%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.62%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.656"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.81%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08308"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.03%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.021"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.70%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.782"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.35%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.85%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-0.57%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9206"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.22%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-0.64%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1958"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.24%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.02%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03881"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'9.37%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'0.99%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001298"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.69%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-4.10%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.441"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.20%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'1.41%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.353"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-4.46%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1376"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.14%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2459"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-9.49%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04404"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.31%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001254"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.38%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004323"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-4.59%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'4.46%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003045"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-23.73%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02790"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'10.05%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05521"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.60%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007777"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.39%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1422"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.67%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008936"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-9.86%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002139"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.01193"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.83%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006955"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.16%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.10%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003180"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'5.12%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.18%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.10%"
$ws.Range("E51").Style = "Normal"
